$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(161, 2).Value = 53925
$ws.Cells.Item(161, 5).Value = 79.37
$ws.Cells.Item(161, 6).Value = 1
$ws.Cells.Item(161, 7).Value = 66.44

$ws.Cells.Item(162, 2).Value = 64350
$ws.Cells.Item(162, 5).Value = 70.63
$ws.Cells.Item(162, 6).Value = 58
$ws.Cells.Item(162, 7).Value = 3853.52

$ws.Cells.Item(163, 2).Value = 57756
$ws.Cells.Item(163, 6).Value = -100
$ws.Cells.Item(163, 7).Value = -6644

$ws.Cells.Item(183, 2).Value = 57552
$ws.Cells.Item(183, 5).Value = 136.86
$ws.Cells.Item(183, 6).Value = -5
$ws.Cells.Item(183, 7).Value = -603.45

$ws.Cells.Item(184, 2).Value = 64329
$ws.Cells.Item(184, 5).Value = 128.32
$ws.Cells.Item(184, 6).Value = 4
$ws.Cells.Item(184, 7).Value = 482.76

$ws.Cells.Item(313, 2).Value = 57854
$ws.Cells.Item(313, 6).Value = 2
$ws.Cells.Item(313, 7).Value = 611.6799999999999

$ws.Cells.Item(314, 2).Value = 62997
$ws.Cells.Item(314, 6).Value = 0
$ws.Cells.Item(314, 7).Value = 0

$ws.Cells.Item(350, 2).Value = 63571
$ws.Cells.Item(350, 5).Value = 152.53
$ws.Cells.Item(350, 6).Value = 19
$ws.Cells.Item(350, 7).Value = 2726.12

$ws.Cells.Item(351, 2).Value = 63531
$ws.Cells.Item(351, 6).Value = 80
$ws.Cells.Item(351, 7).Value = 11478.4

$ws.Cells.Item(352, 2).Value = 57802
$ws.Cells.Item(352, 5).Value = 162.71
$ws.Cells.Item(352, 6).Value = -79
$ws.Cells.Item(352, 7).Value = -11334.92

$ws.Cells.Item(355, 2).Value = 63510
$ws.Cells.Item(355, 5).Value = 50.66
$ws.Cells.Item(355, 6).Value = 159
$ws.Cells.Item(355, 7).Value = 7574.76

$ws.Cells.Item(356, 2).Value = 55356
$ws.Cells.Item(356, 5).Value = 54.04
$ws.Cells.Item(356, 6).Value = -158
$ws.Cells.Item(356, 7).Value = -7527.12

$ws.Cells.Item(372, 2).Value = 57885
$ws.Cells.Item(372, 5).Value = 62.28
$ws.Cells.Item(372, 6).Value = 4
$ws.Cells.Item(372, 7).Value = 208.52

$ws.Cells.Item(373, 2).Value = 63652
$ws.Cells.Item(373, 5).Value = 55.42
$ws.Cells.Item(373, 6).Value = 220
$ws.Cells.Item(373, 7).Value = 11468.6

$ws.Cells.Item(375, 2).Value = 61605
$ws.Cells.Item(375, 5).Value = 133.78
$ws.Cells.Item(375, 6).Value = -13
$ws.Cells.Item(375, 7).Value = -1455.48

$ws.Cells.Item(376, 2).Value = 63563
$ws.Cells.Item(376, 5).Value = 119.04
$ws.Cells.Item(376, 6).Value = 4
$ws.Cells.Item(376, 7).Value = 447.84

$ws.Cells.Item(389, 2).Value = 62865
$ws.Cells.Item(389, 6).Value = 62
$ws.Cells.Item(389, 7).Value = 4948.22

$ws.Cells.Item(390, 2).Value = 57817
$ws.Cells.Item(390, 6).Value = 3
$ws.Cells.Item(390, 7).Value = 239.43

$ws.Cells.Item(400, 2).Value = 57835
$ws.Cells.Item(400, 6).Value = 1
$ws.Cells.Item(400, 7).Value = 59.13

$ws.Cells.Item(401, 2).Value = 62933
$ws.Cells.Item(401, 6).Value = 146
$ws.Cells.Item(401, 7).Value = 8632.98

$ws.Cells.Item(419, 2).Value = 63007
$ws.Cells.Item(419, 6).Value = 886
$ws.Cells.Item(419, 7).Value = 151798.38

$ws.Cells.Item(420, 2).Value = 57856
$ws.Cells.Item(420, 6).Value = 2
$ws.Cells.Item(420, 7).Value = 342.66

$ws.Cells.Item(421, 2).Value = 63008
$ws.Cells.Item(421, 6).Value = 452
$ws.Cells.Item(421, 7).Value = 68328.84

$ws.Cells.Item(422, 2).Value = 57857
$ws.Cells.Item(422, 6).Value = 3
$ws.Cells.Item(422, 7).Value = 453.51

$ws.Cells.Item(536, 2).Value = 47097
$ws.Cells.Item(536, 4).Value = 112.28
$ws.Cells.Item(536, 5).Value = 134.16
$ws.Cells.Item(536, 6).Value = 15
$ws.Cells.Item(536, 7).Value = 1684.2

$ws.Cells.Item(537, 2).Value = 58047
$ws.Cells.Item(537, 4).Value = 105.54
$ws.Cells.Item(537, 5).Value = 126.1
$ws.Cells.Item(537, 6).Value = 47
$ws.Cells.Item(537, 7).Value = 4960.38

$ws.Cells.Item(579, 2).Value = 65069
$ws.Cells.Item(579, 5).Value = 14.3
$ws.Cells.Item(579, 6).Value = 109
$ws.Cells.Item(579, 7).Value = 1466.05

$ws.Cells.Item(580, 2).Value = 53757
$ws.Cells.Item(580, 5).Value = 16.08
$ws.Cells.Item(580, 6).Value = -159
$ws.Cells.Item(580, 7).Value = -2138.55

$ws.Cells.Item(590, 2).Value = 64922
$ws.Cells.Item(590, 5).Value = 20.98
$ws.Cells.Item(590, 6).Value = 174
$ws.Cells.Item(590, 7).Value = 3433.02

$ws.Cells.Item(591, 2).Value = 45706
$ws.Cells.Item(591, 5).Value = 23.58
$ws.Cells.Item(591, 6).Value = -202
$ws.Cells.Item(591, 7).Value = -3985.46

$ws.Cells.Item(593, 2).Value = 45718
$ws.Cells.Item(593, 5).Value = 19.38
$ws.Cells.Item(593, 6).Value = -294
$ws.Cells.Item(593, 7).Value = -4768.68

$ws.Cells.Item(594, 2).Value = 64927
$ws.Cells.Item(594, 5).Value = 17.26
$ws.Cells.Item(594, 6).Value = 282
$ws.Cells.Item(594, 7).Value = 4574.04

$ws.Cells.Item(687, 2).Value = 53319
$ws.Cells.Item(687, 5).Value = 310.64
$ws.Cells.Item(687, 6).Value = -6
$ws.Cells.Item(687, 7).Value = -1643.52

$ws.Cells.Item(688, 2).Value = 64810
$ws.Cells.Item(688, 5).Value = 291.22
$ws.Cells.Item(688, 6).Value = 7
$ws.Cells.Item(688, 7).Value = 1917.44

$ws.Cells.Item(709, 2).Value = 64833
$ws.Cells.Item(709, 5).Value = 34.9
$ws.Cells.Item(709, 6).Value = 98
$ws.Cells.Item(709, 7).Value = 3217.34

$ws.Cells.Item(710, 2).Value = 60025
$ws.Cells.Item(710, 5).Value = 37.22
$ws.Cells.Item(710, 6).Value = -98
$ws.Cells.Item(710, 7).Value = -3217.34

$ws.Cells.Item(720, 2).Value = 64830
$ws.Cells.Item(720, 5).Value = 34.9
$ws.Cells.Item(720, 6).Value = 117
$ws.Cells.Item(720, 7).Value = 3841.11

$ws.Cells.Item(721, 2).Value = 60022
$ws.Cells.Item(721, 5).Value = 37.22
$ws.Cells.Item(721, 6).Value = -113
$ws.Cells.Item(721, 7).Value = -3709.79

$ws.Cells.Item(859, 2).Value = 63150
$ws.Cells.Item(859, 4).Value = 75.68000000000001
$ws.Cells.Item(859, 5).Value = 80.45
$ws.Cells.Item(859, 6).Value = 183
$ws.Cells.Item(859, 7).Value = 13849.44

$ws.Cells.Item(860, 2).Value = 61428
$ws.Cells.Item(860, 4).Value = 69.16
$ws.Cells.Item(860, 5).Value = 73.52
$ws.Cells.Item(860, 6).Value = 1
$ws.Cells.Item(860, 7).Value = 13849.44
